# Fruta / hortaliza, semanal
# Weekly refresh: the rows of data (for each market/category block) get
# reshuffled among each other. Columns D (Fecha), J (Volumen),
# K (Precio mínimo), L (Precio máximo), M (Precio promedio ponderado) and
# P (Precio $/Kg) move between rows 2-14 according to a fixed permutation;
# all other columns (A, B, C, E, F, G, H, I, N, O, Q, R) stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# new row -> source row (where its new data currently lives)
$map = @{
    2  = 6
    3  = 7
    4  = 12
    5  = 14
    6  = 4
    7  = 13
    8  = 3
    9  = 9
    10 = 2
    11 = 8
    12 = 5
    13 = 10
    14 = 11
}

$cols = @("D", "J", "K", "L", "M", "P")

# Snapshot all current values for the columns that move, keyed by row,
# before any writes happen (since source rows and target rows overlap).
$snapshot = @{}
foreach ($r in $map.Keys) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($newRow in $map.Keys) {
    $srcRow = $map[$newRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value2 = $srcVals[$c]
    }
}
